$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 56
$ws.Range("D2").Value = 70
$ws.Range("C3").Value = 62
$ws.Range("D3").Value = 74.5
$ws.Range("C4").Value = 16
$ws.Range("D4").Value = 12
$ws.Range("C5").Value = 222
$ws.Range("D5").Value = 218.5
$ws.Range("C6").Value = 20
$ws.Range("D6").Value = 12.5
$ws.Range("C7").Value = 279
$ws.Range("D7").Value = 323.5
$ws.Range("C8").Value = 23
$ws.Range("D8").Value = 13
$ws.Range("C9").Value = 26
$ws.Range("D9").Value = 26
$ws.Range("C10").Value = 22
$ws.Range("D10").Value = 16
$ws.Range("C11").Value = 8
$ws.Range("D11").Value = 10
$ws.Range("C12").Value = 54
$ws.Range("D12").Value = 59.5
$ws.Range("C13").Value = 223
$ws.Range("D13").Value = 220
$ws.Range("C14").Value = 96
$ws.Range("D14").Value = 88
$ws.Range("C15").Value = 29
$ws.Range("D15").Value = 26
$ws.Range("C16").Value = 32
$ws.Range("D16").Value = 29
$ws.Range("C17").Value = 76
$ws.Range("D17").Value = 71
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 11
$ws.Range("C20").Value = 144
$ws.Range("D20").Value = 152.5
$ws.Range("C21").Value = 49
$ws.Range("D21").Value = 55.5
$ws.Range("C22").Value = 75.25
